# Update "想去人数" (want-to-go count) figures in column F on both the
# "展览" sheet (Worksheets 1) and the mirrored "全部类型" sheet (Worksheets 4).
# Rows G (lowest ticket price) and everything else stay untouched.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 1801
    5  = 777
    13 = 127
    15 = 4198
    19 = 397
    20 = 970
    21 = 1393
    25 = 1941
    28 = 90
}

$sheetIndexes = @(1, 4)

foreach ($sheetIndex in $sheetIndexes) {
    $ws = $wb.Worksheets.Item($sheetIndex)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
